$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 5 data: supriya / supriya@gmail.com / supriya12345 / 18
$ws.Range("A5").Value = "supriya"
$ws.Range("B5").Value = "supriya@gmail.com"
$ws.Range("C5").Value = "supriya12345"
$ws.Range("D5").Value = 18

# Add hyperlink on B5, mirroring the mailto hyperlinks already on B2:B4
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:supriya@gmail.com") | Out-Null

# Make B5 match the Hyperlink cell style used by B2:B4
$ws.Range("B5").Style = $ws.Range("B4").Style

# Move the active selection to C5, matching the saved workbook state
$ws.Range("C5").Select() | Out-Null
